$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# Delete the entire row for program_perc_xpert (row 19), shifting rows 20-25 up.
$ws.Rows.Item(19).Delete()

# The whole-column data validation (smoothness, column C) loses its last row
# when a contained row is deleted; restore it to span the full column again,
# preserving its prompt text.
$c = $ws.Range("C2:C1048575")
$promptTitle = $c.Validation.InputTitle
$promptMsg = $c.Validation.InputMessage
$c.Validation.Delete()
$cFull = $ws.Range("C2:C1048576")
$cFull.Validation.Add(2, 1, 1, "0", "100")
$cFull.Validation.InputTitle = $promptTitle
$cFull.Validation.InputMessage = $promptMsg

$ws.Activate()
$ws.Range("A14").Select()
